$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 52 continues the daily log with the next entry (2025/10/02, 木).
# Force column A to text first so the date-like string "2025/10/02" is
# stored literally (matching the existing rows) instead of being
# auto-converted into a date serial number; then restore the default
# "Normal" style so no stray number-format style id is left behind.
$ws.Cells.Item(52, 1).NumberFormat = "@"
$ws.Cells.Item(52, 1).Value = "2025/10/02"
$ws.Cells.Item(52, 1).Style = "Normal"

$ws.Cells.Item(52, 2).Value = "木"
$ws.Cells.Item(52, 3).Value = 16
$ws.Cells.Item(52, 4).Value = 201
